$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 294; existing rows 294-364 shift down to 295-365
$ws.Rows("294:294").Insert()

# Populate the newly inserted row 294 with its data
$ws.Range("A294").Value = 3
$ws.Range("B294").Value = "Femacal de La Calera"
$ws.Range("C294").Value = "Coquimbo"
$ws.Range("D294").Value = 44722
$ws.Range("E294").Value = 5
$ws.Range("F294").Value = 100114013
$ws.Range("G294").Value = "Zanahoria"
$ws.Range("H294").Value = "Sin especificar"
$ws.Range("I294").Value = "Primera"
$ws.Range("J294").Value = 570
$ws.Range("K294").Value = 6500
$ws.Range("L294").Value = 7000
$ws.Range("M294").Value = 6754
$ws.Range("N294").Value = "$/saco 20 kilos"
$ws.Range("O294").Value = "Chillán"
$ws.Range("P294").Value = 338
$ws.Range("Q294").Value = 20
$ws.Range("R294").Value = "Hortaliza"
